$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Format the header row (row 1): bold, yellow fill, thin box border ---
# Column A keeps the worksheet's existing "text" number format (it already
# carries style index 1, used for SKU values like "0001"), while B:E start
# from the default/general style. Format A1 and B1 individually, then copy
# B1's resulting format onto C1:E1 so every header cell after the first
# reuses the SAME style record instead of the engine minting a fresh one
# per cell.
$a1 = $ws.Cells.Item(1, 1)
$a1.Font.Bold = $true
$a1.Interior.Color = 65535
$a1.Borders.LineStyle = 1

$b1 = $ws.Cells.Item(1, 2)
$b1.Font.Bold = $true
$b1.Interior.Color = 65535
$b1.Borders.LineStyle = 1

$b1.Copy()
$ws.Range("C1:E1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Add the new inventory row ---
$ws.Range("A3").Value = "0002"
$ws.Range("B3").Value = "RFID Wallet Antitheft/Scanning Leather wallet"
$ws.Range("C3").Value = 900
$ws.Range("D3").Value = 1.399
$ws.Range("E3").Value = 1

# --- Restore the selection to where the author left off ---
$ws.Range("C4").Select()
